$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I8").Value = "sv"
$ws.Range("J8").Value = "Statement-opinion"
$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"
$ws.Range("I15").Value = "aa"
$ws.Range("J15").Value = "Agree/Accept"
$ws.Range("I19").Value = "sv"
$ws.Range("J19").Value = "Statement-opinion"
$ws.Range("I20").Value = "%"
$ws.Range("J20").Value = "Uninterpretable"
$ws.Range("I22").Value = "ba"
$ws.Range("J22").Value = "Appreciation"
$ws.Range("I29").Value = "%"
$ws.Range("J29").Value = "Uninterpretable"
$ws.Range("I30").Value = "sd"
$ws.Range("J30").Value = "Statement-non-opinion"
$ws.Range("I50").Value = "ba"
$ws.Range("J50").Value = "Appreciation"
$ws.Range("I52").Value = "sv"
$ws.Range("J52").Value = "Statement-opinion"
$ws.Range("I63").Value = "sd"
$ws.Range("J63").Value = "Statement-non-opinion"
$ws.Range("I71").Value = "b"
$ws.Range("J71").Value = "Acknowledge (Backchannel)"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I95").Value = "sd"
$ws.Range("J95").Value = "Statement-non-opinion"
$ws.Range("I96").Value = "sd"
$ws.Range("J96").Value = "Statement-non-opinion"
$ws.Range("I126").Value = "%"
$ws.Range("J126").Value = "Uninterpretable"
$ws.Range("I132").Value = "sv"
$ws.Range("J132").Value = "Statement-opinion"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I136").Value = "sv"
$ws.Range("J136").Value = "Statement-opinion"
$ws.Range("I146").Value = "%"
$ws.Range("J146").Value = "Uninterpretable"
$ws.Range("I160").Value = "sd"
$ws.Range("J160").Value = "Statement-non-opinion"
$ws.Range("I168").Value = "sv"
$ws.Range("J168").Value = "Statement-opinion"
$ws.Range("I192").Value = "sd"
$ws.Range("J192").Value = "Statement-non-opinion"
$ws.Range("I199").Value = "aa"
$ws.Range("J199").Value = "Agree/Accept"
$ws.Range("I201").Value = "%"
$ws.Range("J201").Value = "Uninterpretable"
$ws.Range("I203").Value = "%"
$ws.Range("J203").Value = "Uninterpretable"
$ws.Range("I204").Value = "%"
$ws.Range("J204").Value = "Uninterpretable"
$ws.Range("I216").Value = "sv"
$ws.Range("J216").Value = "Statement-opinion"
$ws.Range("I223").Value = "ba"
$ws.Range("J223").Value = "Appreciation"
$ws.Range("I225").Value = "%"
$ws.Range("J225").Value = "Uninterpretable"
$ws.Range("I233").Value = "aa"
$ws.Range("J233").Value = "Agree/Accept"
$ws.Range("I235").Value = "aa"
$ws.Range("J235").Value = "Agree/Accept"
$ws.Range("I250").Value = "aa"
$ws.Range("J250").Value = "Agree/Accept"
$ws.Range("I251").Value = "sv"
$ws.Range("J251").Value = "Statement-opinion"
$ws.Range("I252").Value = "b"
$ws.Range("J252").Value = "Acknowledge (Backchannel)"
$ws.Range("I256").Value = "ba"
$ws.Range("J256").Value = "Appreciation"
$ws.Range("I257").Value = "aa"
$ws.Range("J257").Value = "Agree/Accept"
$ws.Range("I267").Value = "sv"
$ws.Range("J267").Value = "Statement-opinion"
$ws.Range("I291").Value = "aa"
$ws.Range("J291").Value = "Agree/Accept"
$ws.Range("I299").Value = "sv"
$ws.Range("J299").Value = "Statement-opinion"
$ws.Range("I303").Value = "ba"
$ws.Range("J303").Value = "Appreciation"
$ws.Range("I306").Value = "sv"
$ws.Range("J306").Value = "Statement-opinion"
$ws.Range("I308").Value = "aa"
$ws.Range("J308").Value = "Agree/Accept"
$ws.Range("I311").Value = "aa"
$ws.Range("J311").Value = "Agree/Accept"
